# Generate Report for Handback
# The 9b2c1858-178e-4bff-b5ee-50b671b78afc.md file has now been handed
# back (it was previously "Ready for handoff"). Update its status to
# "Handed back: in sync with en-US" on every sheet, and record the new
# handback timestamps for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 9b2c1858... file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row 3 is the 9b2c1858... file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("H3").Value = "2016-03-21 22:46:32"

# --- de-de sheet: row 3 is the 9b2c1858... file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("H3").Value = "2016-03-21 22:46:40"
